$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 14 with "Quote History" test data, matching the column layout:
# A=Firstname, B=Lastname, C=VIN, D=Mileage, E=programs, F=Surcharge,
# G=GenerateContract, H=Address, I=Zip Code, J=Phone Number
# Values are written in the same order the shared strings were first
# introduced so the shared string table lines up with the target workbook.
$ws.Range("E14").Value = "New Vehicle - NSC"
$ws.Range("C14").Value = "5FNRL6H27NB019645"
$ws.Range("B14").Value = "Test 14"
$ws.Range("A14").Value = "Automation 14"
$ws.Range("D14").Value = 123
$ws.Range("F14").Value = "one"
$ws.Range("G14").Value = "one"
$ws.Range("H14").Value = "Auto Test"
$ws.Range("I14").Value = 45678
$ws.Range("J14").Value = 9999999

# Match the formatting used by the rest of the data rows (black Calibri
# text, style index reused from the existing rows) on every new cell
# except the VIN column, which keeps the workbook's default style.
$ws.Range("A14:B14").Font.Color = 0
$ws.Range("D14:J14").Font.Color = 0

# Move the active selection to C2, matching the saved cursor position.
[void]$ws.Range("C2").Select()
